{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet habrahabrParagraph = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"habrahabr.ru/company/pentestit/blog/204274\") !== -1) {\n    habrahabrParagraph = items[i];\n    break;\n  }\n}\n\nconst newParagraph = habrahabrParagraph.insertParagraph(\n  \"https://ru.wikipedia.org/wiki/Wireshark\",\n  \"After\"\n);\nawait context.sync();\n\nconst newRange = newParagraph.getRange();\nnewRange.hyperlink = \"https://ru.wikipedia.org/wiki/Wireshark\";\nawait context.sync();\n\nconst ooxmlResult = newParagraph.getOoxml();\nawait context.sync();\nthrow new Error(\"OOXML=\" + JSON.stringify(ooxmlResult.value));\n", "ps1": "$d = $word.ActiveDocument\n$p11 = $d.Paragraphs.Item(11)\n$r = $p11.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n$p12 = $d.Paragraphs.Item(12)\n$p12.Range.LanguageID = \"EnglishUS\"\n$insertPos = $p12.Range.Start\n$d.Range($insertPos, $insertPos).InsertAfter(\"https://ru.wikipedia.org/wiki/Wireshark\")\n\n$p12b = $d.Paragraphs.Item(12)\n$textRange = $p12b.Range\n$textRange.MoveEnd(1, -1)\n$hl = $d.Hyperlinks.Add($textRange, \"https://ru.wikipedia.org/wiki/Wireshark\", \"\", \"\", \"https://ru.wikipedia.org/wiki/Wireshark\")\nthrow \"p12now=[\" + $d.Paragraphs.Item(12).Range.Text + \"]\"\n"}
